# Add a new row (row 5) with 2021 data to Sheet1, mirroring the style of
# existing year rows (A column uses style from A2:A4, i.e. same as the
# header row style which is the "bold/centered/bordered" style index 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the style of A4 (an existing "year" cell) onto A5 so the new label
# cell matches the look of the other year cells in column A.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A5").Value = "2021年"

$ws.Range("B5").Value = 33
$ws.Range("C5").Value = 220449
$ws.Range("D5").Value = 972
$ws.Range("E5").Value = 2020
$ws.Range("F5").Value = 13379
$ws.Range("G5").Value = 33188
$ws.Range("H5").Value = 11177
$ws.Range("I5").Value = 179005
$ws.Range("J5").Value = 39
$ws.Range("K5").Value = 418
$ws.Range("L5").Value = 6560
$ws.Range("M5").Value = 234
